$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Atmos 22 code now working" - adds an alternate wiring configuration (the
# bare/orange/brown wiring scheme) to each of the three 5TM / anemometer
# sections, and renames the section headers to reflect the new hardware
# naming ("Upper Soil (5TM):", "Lower Soil (5TM):",
# "DS2 Anemometer/Atmos 22 Anemometer:").
# ---------------------------------------------------------------------------

# 1) Insert the new rows first (bottom-most block first so the row numbers
#    used for the earlier blocks stay valid).
#    Each of the three blocks gains 4 rows: one italic "OR" divider row plus
#    three new data rows, inserted right before the blank spacer row that
#    follows the block's original 3 data rows.
for ($i = 0; $i -lt 4; $i++) { $ws.Rows.Item(40).Insert() }
for ($i = 0; $i -lt 4; $i++) { $ws.Rows.Item(35).Insert() }

# ---------------------------------------------------------------------------
# Upper Soil (5TM) block - rows 31-38
# ---------------------------------------------------------------------------
$ws.Range("A31").Value = "Upper Soil (5TM):"
$ws.Range("B31").Value = "No intermediate wire necessary. Note: Two sensors present, either wiring configuration possible"

$ws.Range("A32").Value = "Shield"
$ws.Range("B32").Value = "Screw"
$ws.Range("C32").Value = '"U5TM GND" on LHS of LEMS Shield'

$ws.Range("A33").Value = "Red"
$ws.Range("B33").Value = "Screw"
$ws.Range("C33").Value = '"U5TM DATA" on LHS of LEMS Shield'

$ws.Range("A34").Value = "White"
$ws.Range("B34").Value = "Screw"
$ws.Range("C34").Value = '"U5TM POW" on LHS of LEMS Shield'

$ws.Range("A35").Value = "OR"
$ws.Range("A35").Font.Italic = $true

$ws.Range("A36").Value = "Bare/Shield"
$ws.Range("B36").Value = "Screw"
$ws.Range("C36").Value = '"U5TM GND" on LHS of LEMS Shield'

$ws.Range("A37").Value = "Orange"
$ws.Range("B37").Value = "Screw"
$ws.Range("C37").Value = '"U5TM DATA" on LHS of LEMS Shield'

$ws.Range("A38").Value = "Brown"
$ws.Range("B38").Value = "Screw"
$ws.Range("C38").Value = '"U5TM POW" on LHS of LEMS Shield'

$ws.Range("A31").Font.Bold = $true

# ---------------------------------------------------------------------------
# Lower Soil (5TM) block - rows 40-47
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = "Lower Soil (5TM):"
$ws.Range("B40").Value = "No intermediate wire necessary. Note: Two sensors present, either wiring configuration possible"

$ws.Range("A41").Value = "Shield"
$ws.Range("B41").Value = "Screw"
$ws.Range("C41").Value = '"L5TM GND" on RHS of LEMS Shield'

$ws.Range("A42").Value = "Red"
$ws.Range("B42").Value = "Screw"
$ws.Range("C42").Value = '"L5TM DATA" on RHS of LEMS Shield'

$ws.Range("A43").Value = "White"
$ws.Range("B43").Value = "Screw"
$ws.Range("C43").Value = '"L5TM POW" on RHS of LEMS Shield'

$ws.Range("A44").Value = "OR"
$ws.Range("A44").Font.Italic = $true

$ws.Range("A45").Value = "Bare/Shield"
$ws.Range("B45").Value = "Screw"
$ws.Range("C45").Value = '"L5TM GND" on RHS of LEMS Shield'

$ws.Range("A46").Value = "Orange"
$ws.Range("B46").Value = "Screw"
$ws.Range("C46").Value = '"L5TM DATA" on RHS of LEMS Shield'

$ws.Range("A47").Value = "Brown"
$ws.Range("B47").Value = "Screw"
$ws.Range("C47").Value = '"L5TM POW" on RHS of LEMS Shield'

$ws.Range("A40").Font.Bold = $true

# ---------------------------------------------------------------------------
# DS2 Anemometer / Atmos 22 Anemometer block - rows 49-56
# ---------------------------------------------------------------------------
$ws.Range("A49").Value = "DS2 Anemometer/Atmos 22 Anemometer:"
$ws.Range("B49").Value = "No itermediate wire necessary. Note: This anemometer changed in 2018, so there are two names and two different wiring schemes  (as seen below)"

$ws.Range("A50").Value = "Shield"
$ws.Range("B50").Value = "Screw"
$ws.Range("C50").Value = '"GND" on LHS of LEMS Shield'

$ws.Range("A51").Value = "Red"
$ws.Range("B51").Value = "Screw"
$ws.Range("C51").Value = '"SONIC DATA" on LHS of LEMS Shield'

$ws.Range("A52").Value = "White"
$ws.Range("B52").Value = "Screw"
$ws.Range("C52").Value = '"5V" on LHS of LEMS Shield'

$ws.Range("A53").Value = "OR"
$ws.Range("A53").Font.Italic = $true

$ws.Range("A54").Value = "Bare/Shield"
$ws.Range("B54").Value = "Screw"
$ws.Range("C54").Value = '"GND" on LHS of LEMS Shield'

$ws.Range("A55").Value = "Orange"
$ws.Range("B55").Value = "Screw"
$ws.Range("C55").Value = '"SONIC DATA" on LHS of LEMS Shield'

$ws.Range("A56").Value = "Brown"
$ws.Range("B56").Value = "Screw"
$ws.Range("C56").Value = '"5V" on LHS of LEMS Shield'

$ws.Range("A49").Font.Bold = $true

# ---------------------------------------------------------------------------
# Cosmetic sheet-level changes
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 37.1640625

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A41").Select()
